# Update cryptos list values (price + 1h volume/change) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    # Force the cell to keep its text representation even when the
    # new content looks like a number (e.g. "225.66"), matching the
    # original inline-string cells, then restore default formatting.
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.NumberFormat = "General"
}

$ws.Range("D2").Value = "34.079.97"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.784.21"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.59%  "
Set-TextValue $ws.Range("D5") "225.66"
$ws.Range("E5").Value = "  -0.86%  "
Set-TextValue $ws.Range("D6") "0.546"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.46%  "
Set-TextValue $ws.Range("D8") "31.82"
$ws.Range("E8").Value = "  -4.16%  "
Set-TextValue $ws.Range("D9") "0.291"
$ws.Range("E9").Value = "  +1.22%  "
Set-TextValue $ws.Range("D10") "0.0687"
$ws.Range("E10").Value = "  -3.80%  "
Set-TextValue $ws.Range("D11") "0.0942"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "2.045.27"
$ws.Range("E12").Value = "  +0.13%  "
Set-TextValue $ws.Range("D13") "11.16"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "1.800.13"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "34.102.71"
$ws.Range("E15").Value = "  +0.29%  "
Set-TextValue $ws.Range("D16") "0.617"
$ws.Range("E16").Value = "  -1.13%  "
Set-TextValue $ws.Range("D17") "4.18"
$ws.Range("E17").Value = "  +0.98%  "
Set-TextValue $ws.Range("D18") "67.84"
$ws.Range("E18").Value = "  -0.55%  "
Set-TextValue $ws.Range("D19") "245.36"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "0.0₃0776"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  +0.23%  "
Set-TextValue $ws.Range("D22") "10.77"
$ws.Range("E22").Value = "  -0.45%  "
Set-TextValue $ws.Range("D23") "4.08"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("E24").Value = "  -1.91%  "
Set-TextValue $ws.Range("D25") "161.33"
$ws.Range("E25").Value = "  +0.60%  "
Set-TextValue $ws.Range("D26") "7.11"
$ws.Range("E26").Value = "  +0.06%  "
Set-TextValue $ws.Range("D27") "16.26"
$ws.Range("E27").Value = "  -0.69%  "
Set-TextValue $ws.Range("D28") "0.113"
$ws.Range("E28").Value = "  +0.67%  "
Set-TextValue $ws.Range("D29") "1.01"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("E30").Value = "  -0.38%  "
Set-TextValue $ws.Range("D31") "0.0517"
$ws.Range("E31").Value = "  +0.63%  "
Set-TextValue $ws.Range("D32") "3.64"
$ws.Range("E32").Value = "  -0.55%  "
Set-TextValue $ws.Range("D33") "3.58"
$ws.Range("E33").Value = "  +2.05%  "
Set-TextValue $ws.Range("D34") "1.80"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "1.455.67"
$ws.Range("E35").Value = "  +4.14%  "
Set-TextValue $ws.Range("D38") "2.38"
$ws.Range("E38").Value = "  +7.18%  "
$ws.Range("E39").Value = "  -0.67%  "
Set-TextValue $ws.Range("D42") "0.914"
$ws.Range("E42").Value = "  -0.77%  "
Set-TextValue $ws.Range("D43") "2.73"
$ws.Range("E43").Value = "  +1.96%  "
Set-TextValue $ws.Range("D44") "13.30"
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("E45").Value = "  +2.68%  "
Set-TextValue $ws.Range("D46") "6.04"
$ws.Range("E46").Value = "  +3.98%  "
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("E48").Value = "  -0.18%  "
Set-TextValue $ws.Range("D49") "107.20"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").Value = "1.946.12"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  +0.31%  "

# Rows 36/37 swapped rank order: VeChain moves up to rank 36, ImmutableX to rank 37.
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D36") "0.0193"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "0.642"
$ws.Range("E37").Value = "  -1.85%  "

# Rows 40/41 swapped rank order: HuobiToken moves up to rank 40, Aave to rank 41.
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D40") "2.37"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D41") "79.69"
$ws.Range("E41").Value = "  +1.42%  "
